$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark reading progress for 7.17 (D2) and 7.18 (E2) on row 2
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Mark reading progress for 7.20 (G3) on row 3
$ws.Range("G3").Value = 1

# Update the active cell/selection to G9
$ws.Range("G9").Select()
